# Applies the FRConditionLMCDAFHIR mapping update:
#  - Metadata!B9 date bump
#  - "Mapping Table 0" (LM -> CDA): split dateProbleme into dateDebutProbleme/dateFinProbleme
#    (shifts everything below it down one row) + append the trailing "commentaire" row
#  - "Mapping Table 1" (CDA -> FHIR): re-derive from the corrected CDA column + add the two
#    new evidence.detail sub-rows, also appending the trailing "commentaire" row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet: bump the generation Date value
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-01-22T09:24:45+00:00"

# ---------------------------------------------------------------------------
# Mapping Table 0 : FRLMProbleme -> FRCDAProbleme
# ---------------------------------------------------------------------------
$ws0 = $wb.Worksheets.Item("Mapping Table 0")

$table0 = @(
  @("FRLMProbleme", "", "equivalent", "FRCDAProbleme", ""),
  @("FRLMProbleme.identifiant", "", "equivalent", "FRCDAProbleme.id", ""),
  @("FRLMProbleme.type", "", "equivalent", "FRCDAProbleme.code", ""),
  @("FRLMProbleme.description", "", "equivalent", "FRCDAProbleme.text", ""),
  @("FRLMProbleme.problemeObserve", "", "equivalent", "FRCDAProbleme.value", ""),
  @("FRLMProbleme.statut", "", "equivalent", "FRCDAProbleme.statusCode", ""),
  @("FRLMProbleme.dateDebutProbleme", "", "equivalent", "FRCDAProbleme.effectiveTime.low", ""),
  @("FRLMProbleme.dateFinProbleme", "", "equivalent", "FRCDAProbleme.effectiveTime.high", ""),
  @("FRLMProbleme.statutProbleme", "", "equivalent", "FRCDAProbleme.entryRelationship:frStatutDuProbleme", ""),
  @("FRLMProbleme.severite", "", "equivalent", "FRCDAProbleme.entryRelationship:frSeverite", ""),
  @("FRLMProbleme.certitude", "", "equivalent", "FRCDAProbleme.entryRelationship:frCertitude", ""),
  @("FRLMProbleme.statutClinique", "", "equivalent", "FRCDAProbleme.entryRelationship:frStatutCliniqueDuPatient", ""),
  @("FRLMProbleme.reference", "", "equivalent", "FRCDAProbleme.reference", ""),
  @("FRLMProbleme.reference.externalDocument.identifiant", "", "equivalent", "FRCDAProbleme.reference.externalDocument.id", ""),
  @("FRLMProbleme.reference.externalDocument.text.reference", "", "equivalent", "FRCDAProbleme.reference.externalDocument.text.reference", ""),
  @("FRLMProbleme.commentaire", "", "equivalent", "FRCDAProbleme.entryRelationship:frCommentaireER", "")
)

for ($i = 0; $i -lt $table0.Length; $i++) {
  $r = $i + 3
  $row = $table0[$i]
  $ws0.Cells.Item($r, 1).Value = $row[0]
  $ws0.Cells.Item($r, 2).Value = $row[1]
  $ws0.Cells.Item($r, 3).Value = $row[2]
  $ws0.Cells.Item($r, 4).Value = $row[3]
  $ws0.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# Mapping Table 1 : FRCDAProbleme -> FRConditionDocument
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Mapping Table 1")

# Row 2 is the profile-URL banner row (Relationship/Display columns blank)
$ws1.Cells.Item(2, 1).Value = "https://interop.esante.gouv.fr/ig/document/core/StructureDefinition/fr-cda-probleme"
$ws1.Cells.Item(2, 4).Value = "https://interop.esante.gouv.fr/ig/document/core/StructureDefinition/fr-condition-document"

$table1 = @(
  @("FRCDAProbleme", "FRConditionDocument"),
  @("FRCDAProbleme.id", "FRConditionDocument.identifier"),
  @("FRCDAProbleme.code", "FRConditionDocument.category"),
  @("FRCDAProbleme.text", "FRConditionDocument.category.text"),
  @("FRCDAProbleme.value", "FRConditionDocument.code"),
  @("FRCDAProbleme.statusCode", "FRConditionDocument.clinicalStatus"),
  @("FRCDAProbleme.effectiveTime.low", "FRConditionDocument.onsetDateTime"),
  @("FRCDAProbleme.effectiveTime.high", "FRConditionDocument.abatementDateTime"),
  @("FRCDAProbleme.entryRelationship:frStatutDuProbleme", "FRConditionDocument.clinicalStatus"),
  @("FRCDAProbleme.entryRelationship:frSeverite", "FRConditionDocument.severity"),
  @("FRCDAProbleme.entryRelationship:frCertitude", "FRConditionDocument.verificationStatus"),
  @("FRCDAProbleme.entryRelationship:frStatutCliniqueDuPatient", "FRConditionDocument.stage.summary"),
  @("FRCDAProbleme.reference", "FRConditionDocument.evidence.detail"),
  @("FRCDAProbleme.reference.externalDocument.id", "FRConditionDocument.evidence.detail:FRDocumentReferenceDocument.identifier"),
  @("FRCDAProbleme.reference.externalDocument.text.reference", "FRConditionDocument.evidence.detail:FRDocumentReferenceDocument.content.attachment.url"),
  @("FRCDAProbleme.entryRelationship:frCommentaireER", "FRConditionDocument.note")
)

for ($i = 0; $i -lt $table1.Length; $i++) {
  $r = $i + 3
  $row = $table1[$i]
  $ws1.Cells.Item($r, 1).Value = $row[0]
  $ws1.Cells.Item($r, 2).Value = ""
  $ws1.Cells.Item($r, 3).Value = "equivalent"
  $ws1.Cells.Item($r, 4).Value = $row[1]
  $ws1.Cells.Item($r, 5).Value = ""
}
